# Auto-generated edit script: updates Sheets per Halicarnassus_Profits.xlsx diff
# Values are static market-price data (no formulas in this workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 313.22223
$ws.Range("I15").Value = 313.22223
$ws.Range("K15").Value = 939.66669
$ws.Range("M15").Value = -770.66669
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716
$ws.Range("H80").Value = 620.1429000000001
$ws.Range("I80").Value = 588.2
$ws.Range("K80").Value = 1764.6
$ws.Range("M80").Value = -766.6000000000001
$ws.Range("H83").Value = 620.1429000000001
$ws.Range("I83").Value = 588.2
$ws.Range("K83").Value = 5293.8
$ws.Range("M83").Value = -301.8000000000002
$ws.Range("H92").Value = 173
$ws.Range("J92").Value = 80
$ws.Range("L92").Value = 80
$ws.Range("N92").Value = -2576
$ws.Range("H98").Value = 231.25
$ws.Range("I98").Value = 231.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 231.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 1266.75
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 231.25
$ws.Range("I122").Value = 231.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 693.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1756.25
$ws.Range("N122").ClearContents()
$ws.Range("H138").Value = 2741.5715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8498.75
$ws.Range("J61").Value = 8833.333000000001
$ws.Range("L61").Value = 8833.333000000001
$ws.Range("N61").Value = -9257.333000000001
$ws.Range("H63").Value = 8066.1665
$ws.Range("I63").Value = 2099.75
$ws.Range("J63").Value = 19999
$ws.Range("K63").Value = 2099.75
$ws.Range("L63").Value = 19999
$ws.Range("M63").Value = -1413.75
$ws.Range("N63").Value = -21371
$ws.Range("H66").Value = 8066.1665
$ws.Range("I66").Value = 2099.75
$ws.Range("J66").Value = 19999
$ws.Range("K66").Value = 10498.75
$ws.Range("L66").Value = 99995
$ws.Range("M66").Value = -7066.75
$ws.Range("N66").Value = -106859
$ws.Range("H88").Value = 1165
$ws.Range("J88").Value = 1077.625
$ws.Range("L88").Value = 1077.625
$ws.Range("N88").Value = -1889.625
$ws.Range("H91").Value = 1165
$ws.Range("J91").Value = 1077.625
$ws.Range("L91").Value = 1077.625
$ws.Range("N91").Value = -3885.625
$ws.Range("H97").Value = 530.05884
$ws.Range("I97").Value = 576.63635
$ws.Range("K97").Value = 576.63635
$ws.Range("M97").Value = -80.63634999999999
$ws.Range("H102").Value = 3314.3125
$ws.Range("I102").Value = 1569.0834
$ws.Range("K102").Value = 1569.0834
$ws.Range("M102").Value = 52.91660000000002
$ws.Range("H110").Value = 404
$ws.Range("I110").Value = 443.9091
$ws.Range("K110").Value = 443.9091
$ws.Range("M110").Value = 1601.0909
$ws.Range("H122").Value = 3606.2
$ws.Range("I122").Value = 3606.2
$ws.Range("K122").Value = 10818.6
$ws.Range("M122").Value = -8368.599999999999
$ws.Range("H132").Value = 2509.4
$ws.Range("I132").Value = 2509.4
$ws.Range("K132").Value = 7528.200000000001
$ws.Range("M132").Value = -4998.200000000001
$ws.Range("H136").Value = 8498.75
$ws.Range("J136").Value = 8833.333000000001
$ws.Range("L136").Value = 26499.999
$ws.Range("N136").Value = -31599.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1619.625
$ws.Range("I99").Value = 1279.7142
$ws.Range("K99").Value = 1279.7142
$ws.Range("M99").Value = 218.2858000000001
$ws.Range("H105").Value = 1787.1666
$ws.Range("I105").Value = 1823.125
$ws.Range("K105").Value = 1823.125
$ws.Range("M105").Value = -76.125
$ws.Range("H107").Value = 4014.2856
$ws.Range("I107").Value = 2487.2307
$ws.Range("K107").Value = 2487.2307
$ws.Range("M107").Value = -567.2307000000001
$ws.Range("H134").Value = 1884.2858
$ws.Range("I134").Value = 1884.2858
$ws.Range("K134").Value = 5652.857400000001
$ws.Range("M134").Value = -3117.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8600.6
$ws.Range("I62").Value = 8999
$ws.Range("J62").Value = 8501
$ws.Range("K62").Value = 8999
$ws.Range("L62").Value = 8501
$ws.Range("M62").Value = -8375
$ws.Range("N62").Value = -9749
$ws.Range("H65").Value = 8600.6
$ws.Range("I65").Value = 8999
$ws.Range("J65").Value = 8501
$ws.Range("K65").Value = 44995
$ws.Range("L65").Value = 42505
$ws.Range("M65").Value = -41875
$ws.Range("N65").Value = -48745
$ws.Range("H99").Value = 1902.4166
$ws.Range("I99").Value = 1785.8
$ws.Range("K99").Value = 1785.8
$ws.Range("M99").Value = -287.8
$ws.Range("H107").Value = 324.6111
$ws.Range("I107").Value = 221.08333
$ws.Range("K107").Value = 221.08333
$ws.Range("M107").Value = 1698.91667
$ws.Range("H126").Value = 1902.4166
$ws.Range("I126").Value = 1785.8
$ws.Range("K126").Value = 5357.4
$ws.Range("M126").Value = -2887.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 927.5
$ws.Range("I11").Value = 570.3333
$ws.Range("K11").Value = 1710.9999
$ws.Range("M11").Value = -1570.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 343
$ws.Range("I97").Value = 301.07144
$ws.Range("K97").Value = 301.07144
$ws.Range("M97").Value = 194.92856
$ws.Range("H113").Value = 6407.7856
$ws.Range("I113").Value = 5045.6665
$ws.Range("K113").Value = 5045.6665
$ws.Range("M113").Value = -2875.6665
$ws.Range("H132").Value = 1182.3334
$ws.Range("I132").Value = 1182.3334
$ws.Range("K132").Value = 3547.0002
$ws.Range("M132").Value = -1017.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3079
$ws.Range("I7").Value = 2723.75
$ws.Range("K7").Value = 2723.75
$ws.Range("M7").Value = -2611.75
$ws.Range("H40").Value = 2281.7058
$ws.Range("I40").Value = 2281.7058
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2281.7058
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2145.7058
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 5898.25
$ws.Range("J46").Value = 5847.375
$ws.Range("L46").Value = 5847.375
$ws.Range("N46").Value = -6223.375
$ws.Range("H55").Value = 1121.9375
$ws.Range("J55").Value = 1052.4286
$ws.Range("L55").Value = 1052.4286
$ws.Range("N55").Value = -1398.4286
$ws.Range("H63").Value = 44444
$ws.Range("I63").Value = 44444
$ws.Range("K63").Value = 44444
$ws.Range("M63").Value = -43695
$ws.Range("H66").Value = 44444
$ws.Range("I66").Value = 44444
$ws.Range("K66").Value = 133332
$ws.Range("M66").Value = -129588
$ws.Range("H68").Value = 3325
$ws.Range("I68").Value = 1299.5
$ws.Range("J68").Value = 4337.75
$ws.Range("K68").Value = 1299.5
$ws.Range("L68").Value = 4337.75
$ws.Range("M68").Value = -550.5
$ws.Range("N68").Value = -5835.75
$ws.Range("H71").Value = 3325
$ws.Range("I71").Value = 1299.5
$ws.Range("J71").Value = 4337.75
$ws.Range("K71").Value = 6497.5
$ws.Range("L71").Value = 21688.75
$ws.Range("M71").Value = -2753.5
$ws.Range("N71").Value = -29176.75
$ws.Range("H93").Value = 699.8333
$ws.Range("I93").Value = 739.8
$ws.Range("J93").Value = 500
$ws.Range("K93").Value = 739.8
$ws.Range("L93").Value = 500
$ws.Range("M93").Value = 508.2
$ws.Range("N93").Value = -2996
$ws.Range("H100").Value = 6694.8335
$ws.Range("I100").Value = 2667.8
$ws.Range("K100").Value = 2667.8
$ws.Range("M100").Value = -2126.8
$ws.Range("H122").Value = 4060.6667
$ws.Range("I122").Value = 3596
$ws.Range("K122").Value = 10788
$ws.Range("M122").Value = -8338
$ws.Range("H126").Value = 3079
$ws.Range("I126").Value = 2723.75
$ws.Range("K126").Value = 8171.25
$ws.Range("M126").Value = -5701.25
$ws.Range("H136").Value = 2984.625
$ws.Range("I136").Value = 2850.3333
$ws.Range("K136").Value = 8550.999899999999
$ws.Range("M136").Value = -6000.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11250
$ws.Range("I62").Value = 9750
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 9750
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -9126
$ws.Range("N62").Value = -13248
$ws.Range("H65").Value = 11250
$ws.Range("I65").Value = 9750
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 48750
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -45630
$ws.Range("N65").Value = -66240
$ws.Range("H107").Value = 453.92307
$ws.Range("I107").Value = 453.92307
$ws.Range("K107").Value = 1361.76921
$ws.Range("M107").Value = 558.2307900000001
$ws.Range("H136").Value = 2927.742
$ws.Range("I136").Value = 1963.8
$ws.Range("J136").Value = 4680.364
$ws.Range("K136").Value = 5891.4
$ws.Range("L136").Value = 14041.092
$ws.Range("M136").Value = -3341.4
$ws.Range("N136").Value = -19141.092

